$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / date-as-text fields (row 2)
# J2 looks numeric ("001"), so force a text number-format first to keep the
# leading zero / string type intact instead of Excel coercing it to 1.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric fields (row 2)
$ws.Range("O2").Value = 70842023.90000001
$ws.Range("P2").Value = 297.1404314633
$ws.Range("Q2").Value = 289711185.1
$ws.Range("R2").Value = 1215.1672383314
$ws.Range("S2").Value = 88933557.12
$ws.Range("T2").Value = 373.0237234824
$ws.Range("U2").Value = -17009017.03
$ws.Range("V2").Value = -71.342776234
$ws.Range("W2").Value = 493469.88
$ws.Range("X2").Value = 2.0698145675
$ws.Range("Y2").Value = 5512086.91
$ws.Range("Z2").Value = 23.1199476318
$ws.Range("AA2").Value = -30069249.97
$ws.Range("AB2").Value = -126.1227364492
$ws.Range("AC2").Value = 23841260.36
$ws.Range("AD2").Value = 18.8545119368
